# "Your Game Name Here" workload-split spreadsheet:
#   - Title cell (A1) renamed from the placeholder "Your Game Name Here " to
#     the real game title "Going Up!".
#   - The "min" column (D) for a handful of elements that previously had no
#     minimum entered now has one filled in, which ripples through the
#     per-row "Value" formula (column G, =D<n>) and the overall total (G2,
#     =SUM(G4:G23)).
#   - The active selection left on the sheet when it was saved moved to D12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the game.
$ws.Range("A1").Value = "Going Up!"

# Fill in the previously-blank "min" values for several game elements.
$ws.Range("D5").Value = 15
$ws.Range("D6").Value = 5
$ws.Range("D7").Value = 10
$ws.Range("D8").Value = 5
$ws.Range("D11").Value = 5
$ws.Range("D13").Value = 10

# Leave the selection where the author last left it before saving.
$ws.Range("D12").Select()
